$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 811
$ws.Range("D2").Value = "Adult"

$ws.Range("A3").Value = "Nicko"
$ws.Range("B3").Value = "komunist123"
$ws.Range("C3").Value = 1337
$ws.Range("D3").Value = "Minor"

$ws.Range("D3").Select()
